# Add a new worksheet "Sheet2" describing the eth_sb APB-like FSM
# (Current State / Inputs / Next State / Outputs), and adjust the
# view/selection state so Sheet2 becomes the active tab, matching the
# row-height tweak on Sheet1 as well.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")

# --- Add the new sheet right after Sheet1 ---
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Sheet2"

# --- Build 3 reusable style templates (off in a scratch area) by copying
#     existing formats from Sheet1 and tweaking WrapText off, so we don't
#     create stray intermediate cell-style entries. ---
$tHeader = $ws2.Cells.Item(100, 26)   # bold + centered (h+v), no wrap
$ws1.Cells.Item(14, 4).Copy()
$tHeader.PasteSpecial(-4122)
$tHeader.WrapText = $false

$tPlain = $ws2.Cells.Item(101, 26)    # normal font, vertical-centered, no wrap
$ws1.Cells.Item(15, 7).Copy()
$tPlain.PasteSpecial(-4122)
$tPlain.WrapText = $false

$tArial = $ws2.Cells.Item(102, 26)    # Arial Unicode MS 10pt, vertical-centered, no wrap
$ws1.Cells.Item(15, 5).Copy()
$tArial.PasteSpecial(-4122)
$tArial.WrapText = $false

# --- Header row (row 6) ---
$tHeader.Copy()
$ws2.Range("D6:G6").PasteSpecial(-4122)

$headers = @("Current State", "Inputs", "Next State", "Outputs")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws2.Cells.Item(6, 4 + $i).Value = $headers[$i]
}

# --- Body rows 7-13 ---
$rows = @(
    @("IDLE",   "i_eth_sb_psel && fuse_enable",      "SETUP",  "o_eth_sb_pready=0, o_eth_sb_pslverr=0, wr_en=0, rd_en=0"),
    @("IDLE",   "Else",                               "IDLE",   "o_eth_sb_pready=0, o_eth_sb_pslverr=0"),
    @("SETUP",  "!i_eth_sb_psel",                     "IDLE",   "o_eth_sb_pready=0"),
    @("SETUP",  "i_eth_sb_psel",                      "ENABLE", "o_eth_sb_pready=0, wr_en=0, rd_en=0"),
    @("ENABLE", "!i_eth_sb_psel",                     "IDLE",   "Outputs depend on transaction type and FIFO state"),
    @("ENABLE", "i_eth_sb_pwrite && !fifo_full",       "ENABLE", "o_eth_sb_pready=1, wr_en=1, rd_en=0"),
    @("ENABLE", "!i_eth_sb_pwrite && !fifo_empty",     "ENABLE", "o_eth_sb_pready=1, rd_en=1, wr_en=0")
)

# Column style per table column: D=plain, E=arial, F=plain, G=arial
# (except row 8 column E which also uses plain in the source data)
$colStyle = @($tPlain, $tArial, $tPlain, $tArial)

for ($r = 0; $r -lt $rows.Length; $r++) {
    $excelRow = 7 + $r
    $rowData = $rows[$r]
    for ($c = 0; $c -lt $rowData.Length; $c++) {
        $cell = $ws2.Cells.Item($excelRow, 4 + $c)
        $style = $colStyle[$c]
        $style.Copy()
        $cell.PasteSpecial(-4122)
        $cell.Value = $rowData[$c]
    }
}

# Row 8 "Else" cell (E8) uses the plain style rather than Arial in the source.
$tPlain.Copy()
$ws2.Cells.Item(8, 5).PasteSpecial(-4122)
$ws2.Cells.Item(8, 5).Value = "Else"

# Row 11 last column (G11) uses the plain style rather than Arial in the source.
$tPlain.Copy()
$ws2.Cells.Item(11, 7).PasteSpecial(-4122)
$ws2.Cells.Item(11, 7).Value = "Outputs depend on transaction type and FIFO state"

# --- Clean up scratch style template cells ---
$tHeader.Clear()
$tPlain.Clear()
$tArial.Clear()

# --- Column widths on Sheet2 (closest values the engine's pixel-quantized
#     column-width storage can represent to the target 13.28515625 /
#     27.42578125 / 10.42578125 / 51.85546875 character widths) ---
$ws2.Columns.Item(4).ColumnWidth = 12.5
$ws2.Columns.Item(5).ColumnWidth = 26.666666666666668
$ws2.Columns.Item(6).ColumnWidth = 9.666666666666666
$ws2.Columns.Item(7).ColumnWidth = 51.0

# --- Row heights on Sheet1 (30 -> 25.5 for rows 19 and 20) ---
$ws1.Rows.Item(19).RowHeight = 25.5
$ws1.Rows.Item(20).RowHeight = 25.5

# --- Sheet1 selection changes (no longer the tab-selected sheet) ---
$ws1.Range("D28").Select()

# --- Sheet2 becomes the active tab/selection ---
$ws2.Activate()
$ws2.Range("E4").Select()
